# Update loading_percent values for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 9.267912615039105
$ws.Range("C2").Value = 5.176640248673213
$ws.Range("D2").Value = 11.57193575293549
$ws.Range("F2").Value = 29.59414412964679
$ws.Range("G2").Value = 27.57397830717327
$ws.Range("H2").Value = 14.12466470457965
$ws.Range("I2").Value = 20.13318332449921
$ws.Range("J2").Value = 11.32231131363037
$ws.Range("K2").Value = 9.039367830725869
$ws.Range("O2").Value = 21.31730252359417
$ws.Range("B3").Value = 8.897030344599822
$ws.Range("C3").Value = 4.94232216294743
$ws.Range("D3").Value = 11.45614024027583
$ws.Range("F3").Value = 29.66372225482063
$ws.Range("G3").Value = 27.70916398703917
$ws.Range("H3").Value = 14.17883464428974
$ws.Range("I3").Value = 20.24178517478937
$ws.Range("J3").Value = 11.30201281179924
$ws.Range("K3").Value = 8.774221942998841
$ws.Range("O3").Value = 21.41380680445016
$ws.Range("B4").Value = 8.661657169258749
$ws.Range("C4").Value = 4.791831175783424
$ws.Range("D4").Value = 11.38671564810761
$ws.Range("F4").Value = 29.71403485475213
$ws.Range("G4").Value = 27.80123271225656
$ws.Range("H4").Value = 14.21428679866409
$ws.Range("I4").Value = 20.31235662095816
$ws.Range("J4").Value = 11.29178133829331
$ws.Range("K4").Value = 8.6077423286452
$ws.Range("O4").Value = 21.47757430353592
$ws.Range("B5").Value = 8.563948522986021
$ws.Range("C5").Value = 4.728894420923996
$ws.Range("D5").Value = 11.35887456487305
$ws.Range("F5").Value = 29.73644175763586
$ws.Range("G5").Value = 27.8410197342961
$ws.Range("H5").Value = 14.22928530227435
$ws.Range("I5").Value = 20.34209398026746
$ws.Range("J5").Value = 11.28817626280749
$ws.Range("K5").Value = 8.539068866888599
$ws.Range("O5").Value = 21.50469372501991
$ws.Range("B6").Value = 8.547619989852818
$ws.Range("C6").Value = 4.71834829790281
$ws.Range("D6").Value = 11.35427958321235
$ws.Range("F6").Value = 29.74027726006729
$ws.Range("G6").Value = 27.84776301110435
$ws.Range("H6").Value = 14.23180911050833
$ws.Range("I6").Value = 20.3470909973483
$ws.Range("J6").Value = 11.28761180900581
$ws.Range("K6").Value = 8.527618324644351
$ws.Range("O6").Value = 21.50926531619096
$ws.Range("B7").Value = 8.660346507276854
$ws.Range("C7").Value = 4.790988831863076
$ws.Range("D7").Value = 11.38633831383197
$ws.Range("F7").Value = 29.71432933914347
$ws.Range("G7").Value = 27.80176012452224
$ws.Range("H7").Value = 14.21448684037991
$ws.Range("I7").Value = 20.31275370456991
$ws.Range("J7").Value = 11.29173043025188
$ws.Range("K7").Value = 8.606819413235124
$ws.Range("O7").Value = 21.47793545715374
$ws.Range("B8").Value = 9.141694190054244
$ws.Range("C8").Value = 5.097250827439624
$ws.Range("D8").Value = 11.5316795901876
$ws.Range("F8").Value = 29.61655634976348
$ws.Range("G8").Value = 27.61870252503799
$ws.Range("H8").Value = 14.14288785588667
$ws.Range("I8").Value = 20.16982251021598
$ws.Range("J8").Value = 11.31485115764948
$ws.Range("K8").Value = 8.948763028297853
$ws.Range("O8").Value = 21.34963964151654
$ws.Range("B9").Value = 10.01991648261734
$ws.Range("C9").Value = 5.643352640690866
$ws.Range("D9").Value = 11.8286690639849
$ws.Range("F9").Value = 29.48525121296759
$ws.Range("G9").Value = 27.33212758224871
$ws.Range("H9").Value = 14.01985498870229
$ws.Range("I9").Value = 19.9203556177421
$ws.Range("J9").Value = 11.37774453479087
$ws.Range("K9").Value = 9.586369487833348
$ws.Range("O9").Value = 21.13392272687956
$ws.Range("B10").Value = 10.61939238807822
$ws.Range("C10").Value = 6.009228975901968
$ws.Range("D10").Value = 12.05237682084672
$ws.Range("F10").Value = 29.42585213633629
$ws.Range("G10").Value = 27.16635757379004
$ws.Range("H10").Value = 13.94003013879798
$ws.Range("I10").Value = 19.75580472593893
$ws.Range("J10").Value = 11.43441819501323
$ws.Range("K10").Value = 10.03020504897307
$ws.Range("O10").Value = 20.99738331124676
$ws.Range("B11").Value = 10.88122357098071
$ws.Range("C11").Value = 6.167668311525706
$ws.Range("D11").Value = 12.15495713322714
$ws.Range("F11").Value = 29.40691713465492
$ws.Range("G11").Value = 27.10079786101457
$ws.Range("H11").Value = 13.90600572932943
$ws.Range("I11").Value = 19.68500269100382
$ws.Range("J11").Value = 11.46241378904071
$ws.Range("K11").Value = 10.22595336416105
$ws.Range("O11").Value = 20.94005147515902
$ws.Range("B12").Value = 10.97874395801984
$ws.Range("C12").Value = 6.226493298604919
$ws.Range("D12").Value = 12.19388635543989
$ws.Range("F12").Value = 29.40091167562812
$ws.Range("G12").Value = 27.07739872424393
$ws.Range("H12").Value = 13.89345035731811
$ws.Range("I12").Value = 19.65877398308969
$ws.Range("J12").Value = 11.47332770362113
$ws.Range("K12").Value = 10.29913370934545
$ws.Range("O12").Value = 20.91903045328193
$ws.Range("B13").Value = 10.95781461244581
$ws.Range("C13").Value = 6.213876745795934
$ws.Range("D13").Value = 12.18549902229206
$ws.Range("F13").Value = 29.40215322754208
$ws.Range("G13").Value = 27.08237452236062
$ws.Range("H13").Value = 13.8961397536883
$ws.Range("I13").Value = 19.66439690301325
$ws.Range("J13").Value = 11.47096339176812
$ws.Range("K13").Value = 10.28341594786054
$ws.Range("O13").Value = 20.92352701994369
$ws.Range("B14").Value = 10.88927963049858
$ws.Range("C14").Value = 6.172531491230613
$ws.Range("D14").Value = 12.15815836172884
$ws.Range("F14").Value = 29.40639970787303
$ws.Range("G14").Value = 27.09884414401896
$ws.Range("H14").Value = 13.90496619964268
$ws.Range("I14").Value = 19.68283316720322
$ws.Range("J14").Value = 11.46330545353445
$ws.Range("K14").Value = 10.23199322319252
$ws.Range("O14").Value = 20.93830823775652
$ws.Range("B15").Value = 10.84708598989059
$ws.Range("C15").Value = 6.14705303799548
$ws.Range("D15").Value = 12.14142142381315
$ws.Range("F15").Value = 29.40915254276309
$ws.Range("G15").Value = 27.10911839232553
$ws.Range("H15").Value = 13.91041548614567
$ws.Range("I15").Value = 19.69420175948547
$ws.Range("J15").Value = 11.45865527320469
$ws.Range("K15").Value = 10.20037053313954
$ws.Range("O15").Value = 20.94745198300334
$ws.Range("B16").Value = 10.60205645503645
$ws.Range("C16").Value = 5.998711626814954
$ws.Range("D16").Value = 12.04568653340077
$ws.Range("F16").Value = 29.42725248512616
$ws.Range("G16").Value = 27.17084116732672
$ws.Range("H16").Value = 13.9422997497377
$ws.Range("I16").Value = 19.76051333249713
$ws.Range("J16").Value = 11.4326326875164
$ws.Range("K16").Value = 10.01728329298175
$ws.Range("O16").Value = 21.00122643783929
$ws.Range("B17").Value = 10.44890387113469
$ws.Range("C17").Value = 5.905642945989991
$ws.Range("D17").Value = 11.98714145976674
$ws.Range("F17").Value = 29.44042890105851
$ws.Range("G17").Value = 27.21123648612656
$ws.Range("H17").Value = 13.96244571127717
$ws.Range("I17").Value = 19.80223115015094
$ws.Range("J17").Value = 11.41723182445238
$ws.Range("K17").Value = 9.903343946992564
$ws.Range("O17").Value = 21.03544106856765
$ws.Range("B18").Value = 10.35979437265742
$ws.Range("C18").Value = 5.851360382060598
$ws.Range("D18").Value = 11.95354650508173
$ws.Range("F18").Value = 29.44876863227738
$ws.Range("G18").Value = 27.23539732080678
$ws.Range("H18").Value = 13.97424852873804
$ws.Range("I18").Value = 19.82660763813895
$ws.Range("J18").Value = 11.40858246789299
$ws.Range("K18").Value = 9.837232927352947
$ws.Range("O18").Value = 21.0555702825073
$ws.Range("B19").Value = 10.32945037617706
$ws.Range("C19").Value = 5.832852840228424
$ws.Range("D19").Value = 11.94218633889023
$ws.Range("F19").Value = 29.45172295131845
$ws.Range("G19").Value = 27.24373656761605
$ws.Range("H19").Value = 13.97828175800526
$ws.Range("I19").Value = 19.8349266325839
$ws.Range("J19").Value = 11.4056899794072
$ws.Range("K19").Value = 9.814751862302133
$ws.Range("O19").Value = 21.06246289688074
$ws.Range("B20").Value = 10.46531331770706
$ws.Range("C20").Value = 5.915628260687221
$ws.Range("D20").Value = 11.99336578431175
$ws.Range("F20").Value = 29.43894747075988
$ws.Range("G20").Value = 27.20684036718933
$ws.Range("H20").Value = 13.96027884783763
$ws.Range("I20").Value = 19.79775073571162
$ws.Range("J20").Value = 11.41884970233125
$ws.Range("K20").Value = 9.915533067710792
$ws.Range("O20").Value = 21.03175229020007
$ws.Range("B21").Value = 10.90945471454216
$ws.Range("C21").Value = 6.184707574511315
$ws.Range("D21").Value = 12.16618694941091
$ws.Range("F21").Value = 29.40512078781189
$ws.Range("G21").Value = 27.09396780954621
$ws.Range("H21").Value = 13.90236473224128
$ws.Range("I21").Value = 19.67740218333178
$ws.Range("J21").Value = 11.46554634064498
$ws.Range("K21").Value = 10.24712342811535
$ws.Range("O21").Value = 20.93394791351471
$ws.Range("B22").Value = 11.190205403522
$ws.Range("C22").Value = 6.353724353889979
$ws.Range("D22").Value = 12.27961184021028
$ws.Range("F22").Value = 29.38980300427921
$ws.Range("G22").Value = 27.02852072298597
$ws.Range("H22").Value = 13.86643174161419
$ws.Range("I22").Value = 19.60214251167934
$ws.Range("J22").Value = 11.49788431127293
$ws.Range("K22").Value = 10.45830339027296
$ws.Range("O22").Value = 20.87404574192864
$ws.Range("B23").Value = 11.0412540048739
$ws.Range("C23").Value = 6.26414929609652
$ws.Range("D23").Value = 12.21904203377739
$ws.Range("F23").Value = 29.39735661013939
$ws.Range("G23").Value = 27.06268627356803
$ws.Range("H23").Value = 13.88543446405508
$ws.Range("I23").Value = 19.64199945951124
$ws.Range("J23").Value = 11.48046052481718
$ws.Range("K23").Value = 10.3461172931892
$ws.Range("O23").Value = 20.90564833518627
$ws.Range("B24").Value = 10.4578979036015
$ws.Range("C24").Value = 5.911116315712191
$ws.Range("D24").Value = 11.99055156723475
$ws.Range("F24").Value = 29.43961484440059
$ws.Range("G24").Value = 27.20882493457838
$ws.Range("H24").Value = 13.96125779951564
$ws.Range("I24").Value = 19.7997751088716
$ws.Range("J24").Value = 11.41811762136321
$ws.Range("K24").Value = 9.910024248436638
$ws.Range("O24").Value = 21.03341855800053
$ws.Range("B25").Value = 9.790007572851508
$ws.Range("C25").Value = 5.501694321763932
$ws.Range("D25").Value = 11.74722696274833
$ws.Range("F25").Value = 29.51427718951071
$ws.Range("G25").Value = 27.40183552061992
$ws.Range("H25").Value = 14.05128127636149
$ws.Range("I25").Value = 19.98454951848751
$ws.Range("J25").Value = 11.35887346798193
$ws.Range("K25").Value = 9.417901868212164
$ws.Range("O25").Value = 21.18843148777679
